$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("cyan_resolvable_lakes")

# New data rows (order matters for shared-string table ordering)
$ws.Range("A55").Value = "Ross Island Lagoon"
$ws.Range("C55").Value = "Ross Island Lagoon"
$ws.Range("D55").Value = "NO"
$ws.Range("J55").Value = "Ross Island Lagoon"

$ws.Range("A56").Value = "Willamette River (Marquam Brg to St. Johns Brg)"
$ws.Range("C56").Value = "Willamette River (Marquam Brg to St. Johns Brg)"
$ws.Range("D56").Value = "NO"
$ws.Range("J56").Value = "Willamette River (Marquam Brg to St. Johns Brg)"

# New column header
$ws.Range("J1").Value = "NonResolvable"

$ws.Columns.Item(10).EntireColumn.AutoFit()

$ws.Range("A57").Select()
